# Weekly update: a new daily price record for "Arveja Verde" (Macroferia
# Regional de Talca) is inserted as row 17. Excel's row-insert semantics
# shift every existing row from 17 downward by one position (old row 17
# becomes row 18, old row 18 becomes row 19, ..., old row 65 becomes row
# 66), which is exactly what happened in the authored edit - so we let
# Excel do that shift for us instead of rewriting every row by hand.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a fresh blank row at position 17; rows 17-65 move down to 18-66
# and the used range grows to A1:R66 automatically.
$ws.Rows("17:17").Insert()

# Populate the newly inserted row 17 with the new record.
$ws.Range("A17").Value = 5
$ws.Range("B17").Value = "Macroferia Regional de Talca"
$ws.Range("C17").Value = "Maule"
$ws.Range("D17").Value = 44519
$ws.Range("E17").Value = 7
$ws.Range("F17").Value = 100112022
$ws.Range("G17").Value = "Arveja Verde"
$ws.Range("H17").Value = "Sin especificar"
$ws.Range("I17").Value = "Primera"
$ws.Range("J17").Value = 500
$ws.Range("K17").Value = 15000
$ws.Range("L17").Value = 15000
$ws.Range("M17").Value = 15000
$ws.Range("N17").Value = "$/saco 25 kilos"
$ws.Range("O17").Value = "Región del Maule"
$ws.Range("P17").Value = 600
$ws.Range("Q17").Value = 25
$ws.Range("R17").Value = "Hortaliza"
